$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-09-07"

# Update the header label text "2022 (through 09-06)" -> "2022 (through 09-07)"
[void]$ws.Cells.Replace("2022 (through 09-06)", "2022 (through 09-07)")

# Update September 2022 value (row 10) from 31 to 35
$ws.Range("I10").Value = 35

# Update the Total row for the 2022 column (row 14) to reflect new sum
$ws.Range("I14").Value = 1172
